$wb = $excel.ActiveWorkbook

# --- Sheet "isa_template": TEMPLATE metadata -----------------------------
$tpl = $wb.Worksheets.Item("isa_template")

# Version bump
$tpl.Range("B4").Value = "1.0.1"

# Description wording update
$tpl.Range("B5").Value = "Template to describe the nutrients used for plant growth according to MIAPPE. Not all fields may be applicable for every experiment."

# Tags block: split the single "Nutrients" tag into two tags drawn from
# new ontology terms (CHEBI / NCIT) replacing the old MIAPPE tag term.
$tpl.Range("C13").Value = "nutrient"
$tpl.Range("D13").Value = "Growth Medium"

$tpl.Range("C14").Value = "https://bioregistry.io/CHEBI:33284"
$tpl.Range("D14").Value = "https://bioregistry.io/NCIT:C85504"

$tpl.Range("C15").Value = "CHEBI"
$tpl.Range("D15").Value = "NCIT"

# --- Sheet "nutrients": annotation table ---------------------------------
$data = $wb.Worksheets.Item("nutrients")

# Rename "Watering regimen" parameter columns to the updated MIAPPE term
# (regime) and new term numbers (MIAPPE:0161 replacing MIAPPE:0138).
$data.Range("Z1").Value = "Parameter [Watering regime]"
$data.Range("AA1").Value = "Term Source REF (MIAPPE:0161)"
$data.Range("AB1").Value = "Term Accession Number (MIAPPE:0161)"

# Update the UO ontology term accession URLs to the bioregistry.io form.
$data.Range("H2").Value = "https://bioregistry.io/UO:0000309"
$data.Range("V2").Value = "https://bioregistry.io/UO:0000099"
